$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.758.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.91%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.924.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.68%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -1.77%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'335.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.52%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.98%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4675"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.68%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.4143"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.36%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'48.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.50%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.08058"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.03%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.018"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.34%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.30%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.912.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.36%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.024"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.39%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.210"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.99%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'90.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.76%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.9995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.66%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.31%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06589"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.44%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.85%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.9994"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.58%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'29.729.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.53%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.563"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.10%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +7.81%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.199"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.97%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.201.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.96%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'157.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.75%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.150"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.54%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.744"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +7.03%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'117.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.056"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +10.22%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09466"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.65%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.36%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.434"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.00%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.78%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02273"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.71%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'8.480"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.44%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +2.19%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.5920"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.83%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1850"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.79%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.99%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.257"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.71%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.349"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.07520"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.5603"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.21%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'12.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.16%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +2.67%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'112.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.47%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.2996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +11.62%  "
$ws.Range("E51").Style = "Normal"
